$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Developer name
$ws.Range("C3").Value = "Raven Manalastas"

# --- __init__ test rows (7-11): Preconditions=None, Method Inputs=<scenario>, Expected Result=<result>
# NB: column write order below (E, G, F for row 7; E, G for row 8; E,G,F for rows 9-11;
# then F8 deferred to the end of this block) reproduces the shared-string insertion
# order of the authored workbook.
$ws.Range("E7").Value = "None"
$ws.Range("G7").Value = "attributes are set to input values"
$ws.Range("F7").Value = "client_number: 7910`nfirst_name: Jorel`nlast_name: Cruz`nemail_address: jorelcruz@rrc.ca"

$ws.Range("E8").Value = "None"
$ws.Range("G8").Value = "ValueError"

$ws.Range("E9").Value = "None"
$ws.Range("G9").Value = "ValueError"
$ws.Range("F9").Value = "client_number: 7910`nfirst_name: """"`nlast_name: Curz`nemail_address:jorelcruz@rrc.ca"

$ws.Range("E10").Value = "None"
$ws.Range("G10").Value = "ValueError"
$ws.Range("F10").Value = "client_number: 7910`nfirst_name: Jorel`nlast_name: """"`nemail_address: jorelcruz@rrc.ca"

$ws.Range("E11").Value = "None"
$ws.Range("G11").Value = "ValueError"
$ws.Range("F11").Value = "client_number: 7910`nfirst_name: Jorel`nlast_name: Cruz`nemail_address: ""jorelcruzrrc.ca"""

$ws.Range("F8").Value = "client_number: ""INVALID""`nfirst_name: Jorel`nlast_name: Cruz`nemail_address: jorelcruz@rrc.ca"

# --- getter test rows (12-16): Preconditions=<scenario>, Method Inputs=None, Expected Result=<getter value>
$commonPre = "client_number: 7910`nfirst_name: Jorel`nlast_name: Cruz`nemail_address: jorelcruz@rrc.ca"

$ws.Range("E12").Value = $commonPre
$ws.Range("F12").Value = "None"
$ws.Range("G12").Value = 7910
$ws.Range("G12").HorizontalAlignment = -4131

$ws.Range("E13").Value = $commonPre
$ws.Range("F13").Value = "None"
$ws.Range("G13").Value = "Jorel"
$ws.Range("G13").HorizontalAlignment = -4131

$ws.Range("E14").Value = $commonPre
$ws.Range("F14").Value = "None"
$ws.Range("G14").Value = "Cruz"
$ws.Range("G14").HorizontalAlignment = -4131

$ws.Range("E15").Value = $commonPre
$ws.Range("F15").Value = "None"
$ws.Range("G15").HorizontalAlignment = -4131
$ws.Hyperlinks.Add($ws.Range("G15"), "mailto:jorelcruz@rrc.ca", "", "", "jorelcruz@rrc.ca")

$ws.Range("E16").Font.Bold = $true
$ws.Range("F16").Font.Bold = $true
$ws.Range("E16").Value = $commonPre
$ws.Range("F16").Value = "None"
$ws.Range("G16").Font.Bold = $true
$ws.Range("G16").Value = "Cruz, Jorel [7910] - jorelcruz@rrc.ca"
$ws.Range("G16").HorizontalAlignment = -4131

# --- row height touch-ups (mirrors Excel's auto-fit after the content/width edits)
$ws.Rows.Item(2).RowHeight = 73.2
$ws.Rows.Item(12).RowHeight = 66
$ws.Rows.Item(13).RowHeight = 65.4
$ws.Rows.Item(14).RowHeight = 61.2
$ws.Rows.Item(15).RowHeight = 64.8
$ws.Rows.Item(16).RowHeight = 65.4
foreach ($r in 17..28) {
    $ws.Rows.Item($r).RowHeight = 31.2
}

# --- view tweaks
$ws.Application.ActiveWindow.Zoom = 98
$ws.Range("I11").Select()
